$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Status text: "Ready for handoff" -> "Handed back: in sync with en-US"
#    (shown on the Overview sheet for both locales, and on each locale's
#    detail sheet "Status" column)
# ---------------------------------------------------------------------------
$newStatus = "Handed back: in sync with en-US"

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B2").Value = $newStatus
$wsOverview.Range("C2").Value = $newStatus
$wsOverview.Range("B3").Value = $newStatus
$wsOverview.Range("C3").Value = $newStatus

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("C3").Value = $newStatus

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("C3").Value = $newStatus

# ---------------------------------------------------------------------------
# 2. Latest Handback DateTime (column H) now that the handback completed.
#    zh-cn got handed back slightly before de-de, so they end up with
#    different timestamps.
# ---------------------------------------------------------------------------
$wsZhCn.Range("H2").Value = "2016-03-22 12:53:49"
$wsZhCn.Range("H3").Value = "2016-03-22 12:53:49"

$wsDeDe.Range("H2").Value = "2016-03-22 12:53:58"
$wsDeDe.Range("H3").Value = "2016-03-22 12:53:58"

# ---------------------------------------------------------------------------
# 3. Populate "Latest Target File" (F) and "Latest Handback File" (G)
#    hyperlinks for each row on the locale sheets. These mirror the existing
#    Source File Name (A) / Latest Handoff File (D) hyperlinks.
#
#    The existing row-3 hyperlinks (A3/D3) are removed and re-added in
#    between so the final hyperlink list -- and the underlying relationship
#    ids written on save -- stay ordered left-to-right, top-to-bottom:
#    A2, D2, F2, G2, A3, D3, F3, G3.
# ---------------------------------------------------------------------------
function Set-HandbackRow($ws, $row2Md, $row2MdName, $row2Xlf, $row2XlfName, $row3Md, $row3MdName, $row3Xlf, $row3XlfName) {
    $a3Address = $null
    $a3Display = $null
    $d3Address = $null
    $d3Display = $null
    foreach ($hl in $ws.Hyperlinks) {
        $addr = $hl.Range.Address()
        if ($addr -eq '$A$3') {
            $a3Address = $hl.Address
            $a3Display = $hl.TextToDisplay
        }
        if ($addr -eq '$D$3') {
            $d3Address = $hl.Address
            $d3Display = $hl.TextToDisplay
        }
    }

    foreach ($hl in $ws.Hyperlinks) {
        if ($hl.Range.Address() -eq '$A$3') {
            $hl.Delete()
        }
    }
    foreach ($hl in $ws.Hyperlinks) {
        if ($hl.Range.Address() -eq '$D$3') {
            $hl.Delete()
        }
    }

    $ws.Hyperlinks.Add($ws.Range("F2"), $row2Md, "", "", $row2MdName) | Out-Null
    $ws.Hyperlinks.Add($ws.Range("G2"), $row2Xlf, "", "", $row2XlfName) | Out-Null

    $ws.Hyperlinks.Add($ws.Range("A3"), $a3Address, "", "", $a3Display) | Out-Null
    $ws.Hyperlinks.Add($ws.Range("D3"), $d3Address, "", "", $d3Display) | Out-Null

    $ws.Hyperlinks.Add($ws.Range("F3"), $row3Md, "", "", $row3MdName) | Out-Null
    $ws.Hyperlinks.Add($ws.Range("G3"), $row3Xlf, "", "", $row3XlfName) | Out-Null
}

# zh-cn: row 2 (131a7ed5...) and row 3 (388cb62e...)
Set-HandbackRow $wsZhCn `
    "https://github.com/OpenLocalizationTest/oltest/blob/c5624124e8769c4b9f44e6c9b3019166dec92790/e2e/131a7ed5-056e-4060-bc28-41101af30063.md" `
    "131a7ed5-056e-4060-bc28-41101af30063.md" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/b50e4e3b7c04ebda3b1f5b738f321b7ac3c46831/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/131a7ed5-056e-4060-bc28-41101af30063.943fe33fac80d81fc32621ff878cd0d36d500c4f.zh-cn.xlf" `
    "131a7ed5-056e-4060-bc28-41101af30063.943fe33fac80d81fc32621ff878cd0d36d500c4f.zh-cn.xlf" `
    "https://github.com/OpenLocalizationTest/oltest/blob/c5624124e8769c4b9f44e6c9b3019166dec92790/e2e/388cb62e-4b38-4fce-8e47-1c7b6e786d71.md" `
    "388cb62e-4b38-4fce-8e47-1c7b6e786d71.md" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/b50e4e3b7c04ebda3b1f5b738f321b7ac3c46831/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/388cb62e-4b38-4fce-8e47-1c7b6e786d71.e14923f95111d2ddee50c6ec3306a0c813bf63af.zh-cn.xlf" `
    "388cb62e-4b38-4fce-8e47-1c7b6e786d71.e14923f95111d2ddee50c6ec3306a0c813bf63af.zh-cn.xlf"

# de-de: row 2 (131a7ed5...) and row 3 (388cb62e...)
Set-HandbackRow $wsDeDe `
    "https://github.com/OpenLocalizationTest/oltest/blob/c5624124e8769c4b9f44e6c9b3019166dec92790/e2e/131a7ed5-056e-4060-bc28-41101af30063.md" `
    "131a7ed5-056e-4060-bc28-41101af30063.md" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/537ad3fad92a69a7699855009fa84efdede0ec97/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/131a7ed5-056e-4060-bc28-41101af30063.943fe33fac80d81fc32621ff878cd0d36d500c4f.de-de.xlf" `
    "131a7ed5-056e-4060-bc28-41101af30063.943fe33fac80d81fc32621ff878cd0d36d500c4f.de-de.xlf" `
    "https://github.com/OpenLocalizationTest/oltest/blob/c5624124e8769c4b9f44e6c9b3019166dec92790/e2e/388cb62e-4b38-4fce-8e47-1c7b6e786d71.md" `
    "388cb62e-4b38-4fce-8e47-1c7b6e786d71.md" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/537ad3fad92a69a7699855009fa84efdede0ec97/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/388cb62e-4b38-4fce-8e47-1c7b6e786d71.e14923f95111d2ddee50c6ec3306a0c813bf63af.de-de.xlf" `
    "388cb62e-4b38-4fce-8e47-1c7b6e786d71.e14923f95111d2ddee50c6ec3306a0c813bf63af.de-de.xlf"
